# Updated per-team stat block: every league-position row (1-14) gets
# refreshed season totals/rates on both the "per-game" sheet (Sheet1)
# and the "season totals" sheet (Sheet2).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 stat refresh ---
$ws1.Range("B1").Value = 677
$ws1.Range("C1").Value = 216
$ws1.Range("D1").Value = 720
$ws1.Range("E1").Value = 64
$ws1.Range("F1").Value = 0.265
$ws1.Range("G1").Value = 0.796
$ws1.Range("B2").Value = 741
$ws1.Range("C2").Value = 212
$ws1.Range("D2").Value = 699
$ws1.Range("E2").Value = 72
$ws1.Range("G2").Value = 0.786
$ws1.Range("B3").Value = 734
$ws1.Range("C3").Value = 188
$ws1.Range("D3").Value = 677
$ws1.Range("E3").Value = 66
$ws1.Range("F3").Value = 0.245
$ws1.Range("G3").Value = 0.74
$ws1.Range("B4").Value = 631
$ws1.Range("C4").Value = 173
$ws1.Range("D4").Value = 611
$ws1.Range("E4").Value = 87
$ws1.Range("F4").Value = 0.257
$ws1.Range("G4").Value = 0.76
$ws1.Range("B5").Value = 724
$ws1.Range("C5").Value = 208
$ws1.Range("D5").Value = 674
$ws1.Range("E5").Value = 75
$ws1.Range("F5").Value = 0.265
$ws1.Range("G5").Value = 0.814
$ws1.Range("B6").Value = 670
$ws1.Range("C6").Value = 158
$ws1.Range("D6").Value = 587
$ws1.Range("E6").Value = 106
$ws1.Range("F6").Value = 0.268
$ws1.Range("G6").Value = 0.772
$ws1.Range("B7").Value = 720
$ws1.Range("C7").Value = 209
$ws1.Range("D7").Value = 686
$ws1.Range("E7").Value = 91
$ws1.Range("F7").Value = 0.256
$ws1.Range("G7").Value = 0.776
$ws1.Range("B8").Value = 786
$ws1.Range("C8").Value = 232
$ws1.Range("D8").Value = 739
$ws1.Range("E8").Value = 82
$ws1.Range("F8").Value = 0.261
$ws1.Range("G8").Value = 0.809
$ws1.Range("B9").Value = 696
$ws1.Range("C9").Value = 205
$ws1.Range("D9").Value = 659
$ws1.Range("E9").Value = 45
$ws1.Range("F9").Value = 0.264
$ws1.Range("G9").Value = 0.802
$ws1.Range("B10").Value = 736
$ws1.Range("C10").Value = 184
$ws1.Range("D10").Value = 640
$ws1.Range("E10").Value = 77
$ws1.Range("F10").Value = 0.259
$ws1.Range("G10").Value = 0.766
$ws1.Range("B11").Value = 720
$ws1.Range("C11").Value = 241
$ws1.Range("D11").Value = 680
$ws1.Range("E11").Value = 58
$ws1.Range("F11").Value = 0.257
$ws1.Range("G11").Value = 0.809
$ws1.Range("B12").Value = 732
$ws1.Range("C12").Value = 222
$ws1.Range("D12").Value = 723
$ws1.Range("E12").Value = 76
$ws1.Range("F12").Value = 0.261
$ws1.Range("G12").Value = 0.792
$ws1.Range("B13").Value = 728
$ws1.Range("C13").Value = 209
$ws1.Range("D13").Value = 746
$ws1.Range("E13").Value = 99
$ws1.Range("F13").Value = 0.258
$ws1.Range("G13").Value = 0.776
$ws1.Range("B14").Value = 718
$ws1.Range("C14").Value = 254
$ws1.Range("D14").Value = 755
$ws1.Range("E14").Value = 58
$ws1.Range("F14").Value = 0.25

# --- Sheet2 stat refresh ---
$ws2.Range("B1").Value = 89
$ws2.Range("C1").Value = 1557
$ws2.Range("D1").Value = 4.18
$ws2.Range("F1").Value = 112
$ws2.Range("G1").Value = 34
$ws2.Range("B2").Value = 98
$ws2.Range("C2").Value = 1552
$ws2.Range("D2").Value = 4.38
$ws2.Range("E2").Value = 1.26
$ws2.Range("F2").Value = 104
$ws2.Range("G2").Value = 29
$ws2.Range("B3").Value = 78
$ws2.Range("C3").Value = 1286
$ws2.Range("D3").Value = 3.74
$ws2.Range("E3").Value = 1.19
$ws2.Range("F3").Value = 87
$ws2.Range("G3").Value = 43
$ws2.Range("B4").Value = 75
$ws2.Range("C4").Value = 1184
$ws2.Range("D4").Value = 3.34
$ws2.Range("E4").Value = 1.09
$ws2.Range("F4").Value = 66
$ws2.Range("G4").Value = 103
$ws2.Range("B5").Value = 80
$ws2.Range("C5").Value = 1352
$ws2.Range("D5").Value = 4.45
$ws2.Range("E5").Value = 1.29
$ws2.Range("F5").Value = 82
$ws2.Range("G5").Value = 51
$ws2.Range("B6").Value = 73
$ws2.Range("C6").Value = 1357
$ws2.Range("D6").Value = 3.81
$ws2.Range("F6").Value = 84
$ws2.Range("G6").Value = 44
$ws2.Range("B7").Value = 85
$ws2.Range("C7").Value = 1293
$ws2.Range("D7").Value = 4.51
$ws2.Range("E7").Value = 1.3
$ws2.Range("F7").Value = 97
$ws2.Range("G7").Value = -1
$ws2.Range("B8").Value = 89
$ws2.Range("C8").Value = 1554
$ws2.Range("D8").Value = 4.04
$ws2.Range("E8").Value = 1.22
$ws2.Range("F8").Value = 100
$ws2.Range("G8").Value = 35
$ws2.Range("B9").Value = 103
$ws2.Range("C9").Value = 1523
$ws2.Range("D9").Value = 4.05
$ws2.Range("E9").Value = 1.22
$ws2.Range("F9").Value = 110
$ws2.Range("G9").Value = 4
$ws2.Range("B10").Value = 88
$ws2.Range("C10").Value = 1486
$ws2.Range("D10").Value = 3.99
$ws2.Range("F10").Value = 96
$ws2.Range("G10").Value = 3
$ws2.Range("B11").Value = 83
$ws2.Range("C11").Value = 1356
$ws2.Range("D11").Value = 3.86
$ws2.Range("E11").Value = 1.22
$ws2.Range("F11").Value = 81
$ws2.Range("G11").Value = 35
$ws2.Range("B12").Value = 82
$ws2.Range("C12").Value = 1464
$ws2.Range("D12").Value = 3.86
$ws2.Range("E12").Value = 1.21
$ws2.Range("F12").Value = 82
$ws2.Range("G12").Value = 37
$ws2.Range("B13").Value = 70
$ws2.Range("C13").Value = 925
$ws2.Range("D13").Value = 4.42
$ws2.Range("E13").Value = 1.35
$ws2.Range("F13").Value = 61
$ws2.Range("G13").Value = 35
$ws2.Range("B14").Value = 99
$ws2.Range("C14").Value = 1337
$ws2.Range("D14").Value = 3.72
$ws2.Range("E14").Value = 1.19
$ws2.Range("F14").Value = 89
$ws2.Range("G14").Value = 55

# --- Sheet2 view/print changes ---
# Sheet2 stays the active/visible tab; the cursor moves off the full
# A1:G14 selection onto J8 (room reserved for the new sort/print-by-stat feature).
$ws2.Activate()
$ws2.Range("J8").Select()

# Open/touch Page Setup for Sheet2 so a print definition is saved with it,
# matching the sheet now being set up for printing individual stat views.
$ws2.PageSetup.Orientation = 1  # xlPortrait
$ws2.PageSetup.PrintQuality = 300
